$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 198, pushing the existing
# rows 198-204 down to 200-206 (same as the target diff, which keeps the
# original 198-204 block intact three rows further down and adds two new
# "Segunda" quality records dated 2022-01-17 at the top of the block).
$ws.Rows("198:199").Insert()

# New row 198: Melón, Calameño, Segunda
$ws.Cells.Item(198, 1).Value = 4
$ws.Cells.Item(198, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(198, 3).Value = "Los Lagos"
$ws.Cells.Item(198, 4).Value = 44578
$ws.Cells.Item(198, 5).Value = 10
$ws.Cells.Item(198, 6).Value = 100112027
$ws.Cells.Item(198, 7).Value = "Melón"
$ws.Cells.Item(198, 8).Value = "Calameño"
$ws.Cells.Item(198, 9).Value = "Segunda"
$ws.Cells.Item(198, 10).Value = 5000
$ws.Cells.Item(198, 11).Value = 1000
$ws.Cells.Item(198, 12).Value = 1000
$ws.Cells.Item(198, 13).Value = 1000
$ws.Cells.Item(198, 14).Value = "$/unidad"
$ws.Cells.Item(198, 15).Value = "Región del Maule"
$ws.Cells.Item(198, 16).Value = 1000
$ws.Cells.Item(198, 17).Value = 1
$ws.Cells.Item(198, 18).Value = "Hortaliza"

# New row 199: Melón, Tuna, Segunda
$ws.Cells.Item(199, 1).Value = 4
$ws.Cells.Item(199, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(199, 3).Value = "Los Lagos"
$ws.Cells.Item(199, 4).Value = 44578
$ws.Cells.Item(199, 5).Value = 10
$ws.Cells.Item(199, 6).Value = 100112027
$ws.Cells.Item(199, 7).Value = "Melón"
$ws.Cells.Item(199, 8).Value = "Tuna"
$ws.Cells.Item(199, 9).Value = "Segunda"
$ws.Cells.Item(199, 10).Value = 5000
$ws.Cells.Item(199, 11).Value = 1000
$ws.Cells.Item(199, 12).Value = 1000
$ws.Cells.Item(199, 13).Value = 1000
$ws.Cells.Item(199, 14).Value = "$/unidad"
$ws.Cells.Item(199, 15).Value = "Región del Maule"
$ws.Cells.Item(199, 16).Value = 1000
$ws.Cells.Item(199, 17).Value = 1
$ws.Cells.Item(199, 18).Value = "Hortaliza"
